$d = $word.ActiveDocument

# Remove the whole paragraph (including its paragraph mark) that contains
# the stale "Add a way for decision making ..." TODO note - it's been
# superseded now that the bridge has been made in the level design.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Add a way for decision making based on collision detection*") {
        $p.Range.Delete()
        break
    }
}
